# The query field's text run was split by the TokenIteratorFieldRewriterSplit
# parser: the opening "{m:...}" text used to live in a single run, and the
# closing "}" now needs to become its own run, with two new empty runs
# inserted right after the existing _GoBack bookmark (which sits between the
# two pieces of text).

$d = $word.ActiveDocument

# Locate the closing brace of the query expression.
$searchText = "eIDAttribute.name}"
$findRange = $d.Content
$found = $findRange.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the query expression's closing brace"
}
$braceEnd = $findRange.End
$braceStart = $braceEnd - 1

# The _GoBack bookmark currently sits right after the closing brace; move it
# so it sits right before the closing brace instead (i.e. right after
# "...eIDAttribute.name"). Re-adding the bookmark at that collapsed point
# splits the existing run in two without otherwise touching either run's
# text/formatting.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($braceStart, $braceStart))

# Replace the closing "}" character with three runs: two empty runs followed
# by a run containing the closing brace.
$rng = $d.Range($braceStart, $braceEnd)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/word/document.xml" ' + `
       'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
       '<pkg:xmlData>' + `
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
       '<w:body>' + `
       '<w:p>' + `
       '<w:r><w:rPr/><w:t/></w:r>' + `
       '<w:r><w:rPr/><w:t/></w:r>' + `
       '<w:r><w:rPr/><w:t>}</w:t></w:r>' + `
       '</w:p>' + `
       '</w:body></w:document>' + `
       '</pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
